$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "4-Digit 7-Segment Display w/I2C Backpack" row (row 19)
$ws.Rows.Item(19).Delete()
